$wb = $excel.ActiveWorkbook

# --- Sheet "Задача 2" (2nd worksheet) ---
$ws2 = $wb.Worksheets.Item(2)

# B12: was "=H6" (total 7yr compounded growth), now references the annualised
# rate in H7 instead, and is displayed as a percentage.
$ws2.Range("B12").Formula = "=H7"
$ws2.Range("B12").NumberFormat = "0.00%"

# B13: STDEVA range shifts from the raw yearly returns (row 4) to the
# "1+return" row (row 5).
$ws2.Range("B13").Formula = "=STDEVA(B5:H5)"

# B14 (=B12/B13) recalculates automatically from the above changes.

# --- Sheet "Задача 5" (4th worksheet) ---
$ws4 = $wb.Worksheets.Item(4)

# N9 ("resource price" forecast for month 13) updated from 180 to 173.
$ws4.Range("N9").Value = 173.0

# B12:C12 LINEST array formula range changes from B8:M8/B9:M9 to C8:M8/B9:L9
# (regressing against the resource price lagged by one month).
$ws4.Range("B12:C12").FormulaArray = "=LINEST(C8:M8,B9:L9)"

# N8 (=B12*N9+C12) recalculates automatically from the above changes.

$excel.Calculate()
